# Swap values between row 2 and row 3 for columns A, Q, R, S, AW, AX.
# NOTE: cell-by-cell scalar .Value2 assignment (read row2/row3, then write
# back individually) triggers a shared-string write-back bug in this
# runtime's save path when multiple cells end up referencing pre-existing
# duplicate string values. Swapping each column's 2-cell block via a single
# array assignment avoids that bug and reliably applies the correct swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = "A", "Q", "R", "S", "AW", "AX"

foreach ($col in $cols) {
    $rng = $ws.Range("$col`2:$col`3")
    $arr = $rng.Value2
    $newArr = New-Object 'object[,]' 2,1
    $newArr[0,0] = $arr[2,1]
    $newArr[1,0] = $arr[1,1]
    $rng.Value2 = $newArr
}
